$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.094.94'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '3.394.52'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''573.31'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '''142.16'
$ws.Range("D7").Value = '3.394.70'
$ws.Range("E7").Value = '  -2.04%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.475'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("E11").Value = '  -2.95%  '
$ws.Range("D12").Value = '''0.394'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '3.976.15'
$ws.Range("E13").Value = '  -2.05%  '
$ws.Range("E14").Value = '  +2.39%  '
$ws.Range("D15").Value = '''27.80'
$ws.Range("E15").Value = '  -3.91%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000171'
$ws.Range("E16").Value = '  -2.34%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.370.50'
$ws.Range("E17").Value = '  -2.45%  '
$ws.Range("D18").Value = '61.152.40'
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("D19").Value = '''6.12'
$ws.Range("E19").Value = '  -4.58%  '
$ws.Range("D20").Value = '''13.77'
$ws.Range("E20").Value = '  -4.82%  '
$ws.Range("D21").Value = '''8.91'
$ws.Range("E21").Value = '  -6.15%  '
$ws.Range("D22").Value = '''382.96'
$ws.Range("E22").Value = '  -4.22%  '
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("D24").Value = '''74.66'
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '''0.0000117'
$ws.Range("E26").Value = '  -5.68%  '
$ws.Range("D27").Value = '3.533.67'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '''7.35'
$ws.Range("E30").Value = '  -4.56%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''2.16'
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''7.98'
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("E33").Value = '  -6.05%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '''23.44'
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("D36").Value = '''6.98'
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("D37").Value = '''166.46'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").Value = '3.426.94'
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("D39").Value = '''5.00'
$ws.Range("E39").Value = '  -3.83%  '
$ws.Range("D40").Value = '''1.48'
$ws.Range("E40").Value = '  -5.42%  '
$ws.Range("D41").Value = '''0.0770'
$ws.Range("E41").Value = '  -3.22%  '
$ws.Range("D42").Value = '''26.96'
$ws.Range("E42").Value = '  -5.96%  '
$ws.Range("D43").Value = '''0.781'
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").Value = '''41.90'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").Value = '''4.41'
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("E48").Value = '  -1.87%  '
$ws.Range("D49").Value = '2.475.79'
$ws.Range("E49").Value = '  -6.59%  '
$ws.Range("D50").Value = '''6.79'
$ws.Range("E50").Value = '  -2.78%  '
$ws.Range("D51").Value = '''22.92'
$ws.Range("E51").Value = '  -0.66%  '
